$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as exact text in the source data;
# force text format so Excel does not silently coerce them to floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.87"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.97"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.199"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05744"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.449"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.229"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8132"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8687"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1376"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06938"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03028"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09328"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.817"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001523"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04719"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005967"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006203"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001235"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008694"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.583"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03714"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006236"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1049"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002598"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008368"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4298"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002058"

# Coin / link / volume-label text columns
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
